$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 392.75
$ws.Range("I33").Value = 434.7143
$ws.Range("J33").Value = 334
$ws.Range("K33").Value = 434.7143
$ws.Range("L33").Value = 334
$ws.Range("M33").Value = -205.7143
$ws.Range("N33").Value = -792
$ws.Range("H112").Value = 58826292
$ws.Range("I112").Value = 333334100
$ws.Range("J112").Value = 3186.4285
$ws.Range("K112").Value = 1000002300
$ws.Range("L112").Value = 9559.2855
$ws.Range("M112").Value = -1000001192
$ws.Range("N112").Value = -11775.2855
$ws.Range("H114").Value = 39721.4
$ws.Range("J114").Value = 39721.4
$ws.Range("L114").Value = 39721.4
$ws.Range("N114").Value = -48399.4
$ws.Range("H135").Value = 53639.08
$ws.Range("I135").Value = 29634.656
$ws.Range("J135").Value = 333690.66
$ws.Range("K135").Value = 266711.904
$ws.Range("L135").Value = 3003215.94
$ws.Range("M135").Value = -264176.904
$ws.Range("N135").Value = -3008285.94
$ws.Range("H138").Value = 3637901
$ws.Range("I138").Value = 1253.25
$ws.Range("J138").Value = 5715985.5
$ws.Range("K138").Value = 3759.75
$ws.Range("L138").Value = 17147956.5
$ws.Range("M138").Value = 1380.25
$ws.Range("N138").Value = -17158236.5

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 7593.94
$ws.Range("I32").Value = 4373.5664
$ws.Range("J32").Value = 23316.941
$ws.Range("K32").Value = 4373.5664
$ws.Range("L32").Value = 23316.941
$ws.Range("M32").Value = -4086.5664
$ws.Range("N32").Value = -23890.941
$ws.Range("H70").Value = 0
$ws.Range("J70").Value = 0
$ws.Range("L70").Value = 0
$ws.Range("N70").ClearContents()
$ws.Range("H73").Value = 0
$ws.Range("J73").Value = 0
$ws.Range("L73").Value = 0
$ws.Range("N73").ClearContents()
$ws.Range("H135").Value = 45539.715
$ws.Range("J135").Value = 45539.715
$ws.Range("L135").Value = 45539.715
$ws.Range("N135").Value = -55679.715

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H123").Value = 19140.643
$ws.Range("J123").Value = 19140.643
$ws.Range("L123").Value = 19140.643
$ws.Range("N123").Value = -28940.643
$ws.Range("H134").Value = 1939.283
$ws.Range("I134").Value = 1851.6818
$ws.Range("J134").Value = 2367.5557
$ws.Range("K134").Value = 5555.0454
$ws.Range("L134").Value = 7102.6671
$ws.Range("M134").Value = -3020.0454
$ws.Range("N134").Value = -12172.6671

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H86").Value = 2963.5
$ws.Range("I86").Value = 2900
$ws.Range("J86").Value = 2972.5715
$ws.Range("K86").Value = 2900
$ws.Range("L86").Value = 2972.5715
$ws.Range("M86").Value = -1777
$ws.Range("N86").Value = -5218.5715
$ws.Range("H89").Value = 2963.5
$ws.Range("I89").Value = 2900
$ws.Range("J89").Value = 2972.5715
$ws.Range("K89").Value = 14500
$ws.Range("L89").Value = 14862.8575
$ws.Range("M89").Value = -8884
$ws.Range("N89").Value = -26094.8575
$ws.Range("H132").Value = 26074.428
$ws.Range("I132").Value = 1943.3667
$ws.Range("J132").Value = 86402.086
$ws.Range("K132").Value = 5830.1001
$ws.Range("L132").Value = 259206.258
$ws.Range("M132").Value = -3300.1001
$ws.Range("N132").Value = -264266.258
$ws.Range("H134").Value = 21362.873
$ws.Range("I134").Value = 1756.3243
$ws.Range("J134").Value = 61665.223
$ws.Range("K134").Value = 5268.9729
$ws.Range("L134").Value = 184995.669
$ws.Range("M134").Value = -2733.9729
$ws.Range("N134").Value = -190065.669

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 178.26666
$ws.Range("J2").Value = 73.71429000000001
$ws.Range("L2").Value = 442.28574
$ws.Range("N2").Value = -668.28574
$ws.Range("H32").Value = 730
$ws.Range("J32").Value = 900
$ws.Range("L32").Value = 2700
$ws.Range("N32").Value = -3266
$ws.Range("H46").Value = 1646.6666
$ws.Range("J46").Value = 1920
$ws.Range("L46").Value = 5760
$ws.Range("N46").Value = -5942
$ws.Range("H118").Value = 3316.889
$ws.Range("I118").Value = 450.66666
$ws.Range("J118").Value = 4750
$ws.Range("K118").Value = 1351.99998
$ws.Range("L118").Value = 14250
$ws.Range("M118").Value = -108.9999800000001
$ws.Range("N118").Value = -16736
$ws.Range("H131").Value = 830.93243
$ws.Range("J131").Value = 846.5571
$ws.Range("L131").Value = 2539.6713
$ws.Range("N131").Value = -12619.6713

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 0
$ws.Range("K4").Value = 0
$ws.Range("M4").ClearContents()
$ws.Range("H70").Value = 43725.926
$ws.Range("I70").Value = 96727.27
$ws.Range("J70").Value = 7287.5
$ws.Range("K70").Value = 96727.27
$ws.Range("L70").Value = 7287.5
$ws.Range("M70").Value = -96457.27
$ws.Range("N70").Value = -7827.5
$ws.Range("H73").Value = 43725.926
$ws.Range("I73").Value = 96727.27
$ws.Range("J73").Value = 7287.5
$ws.Range("K73").Value = 96727.27
$ws.Range("L73").Value = 7287.5
$ws.Range("M73").Value = -95791.27
$ws.Range("N73").Value = -9159.5
$ws.Range("H80").Value = 4063.077
$ws.Range("I80").Value = 3380
$ws.Range("J80").Value = 4490
$ws.Range("K80").Value = 3380
$ws.Range("L80").Value = 4490
$ws.Range("M80").Value = -2382
$ws.Range("N80").Value = -6486
$ws.Range("H83").Value = 4063.077
$ws.Range("I83").Value = 3380
$ws.Range("J83").Value = 4490
$ws.Range("K83").Value = 16900
$ws.Range("L83").Value = 22450
$ws.Range("M83").Value = -11908
$ws.Range("N83").Value = -32434

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 9262133
$ws.Range("I7").Value = 13891375
$ws.Range("J7").Value = 3648.7778
$ws.Range("K7").Value = 13891375
$ws.Range("L7").Value = 3648.7778
$ws.Range("M7").Value = -13891263
$ws.Range("N7").Value = -3872.7778
$ws.Range("H122").Value = 3602.195
$ws.Range("I122").Value = 3303.6365
$ws.Range("J122").Value = 3711.6667
$ws.Range("K122").Value = 9910.9095
$ws.Range("L122").Value = 11135.0001
$ws.Range("M122").Value = -7460.9095
$ws.Range("N122").Value = -16035.0001
$ws.Range("H126").Value = 9262133
$ws.Range("I126").Value = 13891375
$ws.Range("J126").Value = 3648.7778
$ws.Range("K126").Value = 41674125
$ws.Range("L126").Value = 10946.3334
$ws.Range("M126").Value = -41671655
$ws.Range("N126").Value = -15886.3334
$ws.Range("H132").Value = 68674.8
$ws.Range("I132").Value = 1168
$ws.Range("K132").Value = 3504
$ws.Range("M132").Value = -974

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 2350.4
$ws.Range("I126").Value = 2072
$ws.Range("K126").Value = 6216
$ws.Range("M126").Value = -3746
$ws.Range("H132").Value = 28483.562
$ws.Range("I132").Value = 18018.275
$ws.Range("J132").Value = 68949.336
$ws.Range("K132").Value = 54054.825
$ws.Range("L132").Value = 206848.008
$ws.Range("M132").Value = -51524.825
$ws.Range("N132").Value = -211908.008
